$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("LC","LD","LE","LF","LG","LH")

$data = @{
  1 = @(11000,11010,11018,11027,11036,11041)
  2 = @(2023,2023,2023,2023,2023,2023)
  3 = @(4,5,6,7,8,9)
  4 = @(0,1,1,0,1,0)
  5 = @(1,0,0,0,1,1)
  6 = @(63,63,68,95,103,39)
  7 = @(73,66,105,71,49,87)
  8 = @(-10,-3,-37,24,54,-48)
  9 = @(0,0,0,1,1,0)
  10 = @(3,13,17,18,14,16)
  11 = @(207,177,194,206,238,194)
  12 = @(164,159,136,179,166,175)
  13 = @(371,336,330,385,404,369)
  14 = @(1.26,1.11,1.43,1.15,1.43,1.11)
  15 = @(90,66,94,91,130,97)
  16 = @(66,53,49,53,60,42)
  17 = @(36,25,52,47,26,13)
  18 = @(16,19,17,14,20,24)
  19 = @(12,23,16,18,11,17)
  20 = @(9,9,10,14,15,4)
  21 = @(3,7,9,11,9,4)
  22 = @(8,5,6,11,9,14)
  23 = @(1,4,2,0,4,1)
  24 = @(18,18,18,25,28,19)
  25 = @(50,50,55.6,56,53.6,21.1)
  26 = @(41.22,37.33,33,27.5,26.93,92.25)
  27 = @(20.61,18.67,18.33,15.4,14.43,19.42)
  28 = @(44,40,42,47,36,28)
  29 = @(50,59,41,46,47,53)
  30 = @(38,47,32,40,32,37)
  31 = @(45,46,40,52,59,46)
  32 = @(2.5,2.56,2.22,2.08,2.11,2.42)
  33 = @(5,5.11,4,3.71,3.93,11.5)
  34 = @(37.8,30.4,40,48.1,40.7,39.1)
  35 = @(20,19.6,25,26.9,25.4,8.7)
  36 = @(188.7,189,189.7,189.1,189.2,188.7)
  37 = @(87.3,87.3,88.3,88.4,88.1,87.1)
  38 = @(25.41,25.33,25,24.66,24.66,24.8)
  39 = @(90.4,89.6,87,84.1,86.4,86.1)
  40 = @(7,8,8,8,8,8)
  41 = @(7,6,7,7,7,8)
  42 = @(5,5,4,4,4,3)
  43 = @(4,4,4,4,4,4)
  44 = @(130,138,108,136,111,110)
  45 = @(235,195,212,246,280,255)
  46 = @(283,244,258,301,312,281)
  47 = @(76.3,72.6,78.2,78.2,77.2,76.2)
  48 = @(50,59,41,46,47,53)
  49 = @(6,7,7,9,6,6)
  50 = @(7,6,9,15,18,12)
  51 = @(44,40,42,47,36,28)
  52 = @(38,47,32,40,32,37)
  53 = @(43,42,33,39,29,36)
  54 = @(7,3,9,14,10,1)
  55 = @(3,7,9,11,9,4)
  56 = @(33.3,77.8,90,78.6,60,100)
  57 = @(219,212,204,203,191,253)
  58 = @(134,126,113,159,131,124)
  59 = @(353,338,317,362,322,377)
  60 = @(1.63,1.68,1.81,1.28,1.46,2.04)
  61 = @(105,93,101,99,97,128)
  62 = @(80,73,48,60,40,52)
  63 = @(39,67,24,34,34,41)
  64 = @(12,23,16,18,11,17)
  65 = @(16,19,17,14,20,24)
  66 = @(10,9,16,10,6,13)
  67 = @(5,6,7,7,5,11)
  68 = @(12,11,9,10,10,8)
  69 = @(1,1,0,1,3,1)
  70 = @(23,21,25,21,19,22)
  71 = @(43.5,42.9,64,47.6,31.6,59.1)
  72 = @(35.3,37.56,19.81,36.2,53.67,29)
  73 = @(15.35,16.1,12.68,17.24,16.95,17.14)
  74 = @(27,31,36,30,16,27)
  75 = @(49,49,38,45,66,55)
  76 = @(36,36,29,38,43,42)
  77 = @(48,56,48,50,38,51)
  78 = @(2.09,2.67,1.92,2.38,2,2.32)
  79 = @(4.8,6.22,3,5,6.33,3.92)
  80 = @(45.8,35.7,52.1,40,42.1,41.2)
  81 = @(20.8,16.1,33.3,20,15.8,25.5)
  82 = @(187.5,188.2,189.7,189,187.2,187.2)
  83 = @(86.3,89.2,87.6,86.6,85.1,86.8)
  84 = @(26.8,25.58,25.49,26.66,26.41,26.33)
  85 = @(92.1,96,99.3,116.1,92.2,117.5)
  86 = @(6,5,8,7,12,2)
  87 = @(5,9,5,2,1,10)
  88 = @(8,5,4,5,4,3)
  89 = @(4,4,6,9,6,8)
  90 = @(136,127,132,124,98,115)
  91 = @(215,200,192,232,219,246)
  92 = @(254,255,240,270,245,294)
  93 = @(72,75.4,75.7,74.6,76.1,78)
  94 = @(49,49,38,45,66,55)
  95 = @(16,13,17,8,6,8)
  96 = @(8,12,14,9,9,11)
  97 = @(27,31,36,30,16,27)
  98 = @(36,36,29,38,43,42)
  99 = @(41,45,40,45,41,31)
  100 = @(6,8,5,4,13,6)
  101 = @(5,6,7,7,5,11)
  102 = @(50,66.7,43.8,70,83.3,84.6)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = $cols[$i] + $row
    $ws.Range($addr).Value2 = $vals[$i]
  }
}

"done"